$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Done?" status (column K) for cards 21-32 (rows 23-34).
# The SleightOfHand card (row 31, card id 29) has now been implemented ("Y"),
# while the rest of this newly tracked range is not yet done ("N").
$ws.Range("K23").Value = "N"
$ws.Range("K24").Value = "N"
$ws.Range("K25").Value = "N"
$ws.Range("K26").Value = "N"
$ws.Range("K27").Value = "N"
$ws.Range("K28").Value = "N"
$ws.Range("K29").Value = "N"
$ws.Range("K30").Value = "N"
$ws.Range("K31").Value = "Y"
$ws.Range("K32").Value = "N"
$ws.Range("K33").Value = "N"
$ws.Range("K34").Value = "N"

# Update the active selection to K31 to match the edited cell.
$ws.Range("K31").Select()
